$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 102.4929963333333
$ws.Range("H2").Value = 307.478989
$ws.Range("I2").Value = 0.2065071987599813
$ws.Range("J2").Value = 0.2065071987599814
$ws.Range("M2").Value = 61.04160633333334
$ws.Range("N2").Value = 183.124819
$ws.Range("O2").Value = 0.2043613460574534
$ws.Range("P2").Value = 0.2043613460574534
$ws.Range("Q2").Value = 6256.337134103112
$ws.Range("R2").Value = 56307.034206928
$ws.Range("S2").Value = 0.04220208910914387
$ws.Range("T2").Value = 0.04220208910914387
$ws.Range("G3").Value = 102.4929963333333
$ws.Range("H3").Value = 307.478989
$ws.Range("I3").Value = 0.2065071987599813
$ws.Range("J3").Value = 0.2065071987599814
$ws.Range("O3").Value = 0.3559304658284363
$ws.Range("P3").Value = 0.3559304658284363
$ws.Range("Q3").Value = 10896.48817391829
$ws.Range("R3").Value = 98068.39356526463
$ws.Range("S3").Value = 0.07350220345156563
$ws.Range("T3").Value = 0.07350220345156566
$ws.Range("G4").Value = 102.4929963333333
$ws.Range("H4").Value = 307.478989
$ws.Range("I4").Value = 0.2065071987599813
$ws.Range("J4").Value = 0.2065071987599814
$ws.Range("M4").Value = 131.3384093333333
$ws.Range("N4").Value = 394.015228
$ws.Range("O4").Value = 0.4397081881141102
$ws.Range("P4").Value = 0.4397081881141103
$ws.Range("Q4").Value = 13461.26710622716
$ws.Range("R4").Value = 121151.4039560445
$ws.Range("S4").Value = 0.09080290619927182
$ws.Range("T4").Value = 0.09080290619927184
$ws.Range("I5").Value = 0.581825957350084
$ws.Range("J5").Value = 0.5818259573500841
$ws.Range("M5").Value = 61.04160633333334
$ws.Range("N5").Value = 183.124819
$ws.Range("O5").Value = 0.2043613460574534
$ws.Range("P5").Value = 0.2043613460574534
$ws.Range("Q5").Value = 17626.98523059833
$ws.Range("R5").Value = 158642.8670753849
$ws.Range("S5").Value = 0.1189027358152297
$ws.Range("T5").Value = 0.1189027358152297
$ws.Range("I6").Value = 0.581825957350084
$ws.Range("J6").Value = 0.5818259573500841
$ws.Range("O6").Value = 0.3559304658284363
$ws.Range("P6").Value = 0.3559304658284363
$ws.Range("S6").Value = 0.2070895840306913
$ws.Range("T6").Value = 0.2070895840306914
$ws.Range("I7").Value = 0.581825957350084
$ws.Range("J7").Value = 0.5818259573500841
$ws.Range("M7").Value = 131.3384093333333
$ws.Range("N7").Value = 394.015228
$ws.Range("O7").Value = 0.4397081881141102
$ws.Range("P7").Value = 0.4397081881141103
$ws.Range("Q7").Value = 37926.59368901181
$ws.Range("R7").Value = 341339.3432011063
$ws.Range("S7").Value = 0.255833637504163
$ws.Range("T7").Value = 0.2558336375041631
$ws.Range("G8").Value = 105.053815
$ws.Range("H8").Value = 315.161445
$ws.Range("I8").Value = 0.2116668438899346
$ws.Range("J8").Value = 0.2116668438899346
$ws.Range("M8").Value = 61.04160633333334
$ws.Range("N8").Value = 183.124819
$ws.Range("O8").Value = 0.2043613460574534
$ws.Range("P8").Value = 0.2043613460574534
$ws.Range("Q8").Value = 6412.65361904483
$ws.Range("R8").Value = 57713.88257140347
$ws.Range("S8").Value = 0.0432565211330799
$ws.Range("T8").Value = 0.04325652113307991
$ws.Range("G9").Value = 105.053815
$ws.Range("H9").Value = 315.161445
$ws.Range("I9").Value = 0.2116668438899346
$ws.Range("J9").Value = 0.2116668438899346
$ws.Range("O9").Value = 0.3559304658284363
$ws.Range("P9").Value = 0.3559304658284363
$ws.Range("Q9").Value = 11168.74024298779
$ws.Range("R9").Value = 100518.6621868901
$ws.Range("S9").Value = 0.07533867834617933
$ws.Range("T9").Value = 0.07533867834617936
$ws.Range("G10").Value = 105.053815
$ws.Range("H10").Value = 315.161445
$ws.Range("I10").Value = 0.2116668438899346
$ws.Range("J10").Value = 0.2116668438899346
$ws.Range("M10").Value = 131.3384093333333
$ws.Range("N10").Value = 394.015228
$ws.Range("O10").Value = 0.4397081881141102
$ws.Range("P10").Value = 0.4397081881141103
$ws.Range("Q10").Value = 13797.60095649827
$ws.Range("R10").Value = 124178.4086084845
$ws.Range("S10").Value = 0.09307164441067538
$ws.Range("T10").Value = 0.09307164441067539
